$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New columns F and G: set widths so <col> entries get created, same
#    width (18 chars) as the existing data columns B:E.
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 17.17
$ws.Columns.Item(7).ColumnWidth = 17.17

# ---------------------------------------------------------------------
# 2. Header row timestamps for the two new runs (same style as E1).
# ---------------------------------------------------------------------
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("F1").Value = "2025-12-14 19:36:19"
$ws.Range("G1").Value = "2025-12-14 19:41:18"

# ---------------------------------------------------------------------
# 3. Seed the "blank/no-result" style (text number format, no fill) on
#    F2, matching the run where column F didn't have data for that
#    model yet. This creates the new shared cellXf used throughout.
# ---------------------------------------------------------------------
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Copy()
$ws.Range("F3:F5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. Results for runs F (2025-12-14 19:36:19) and G (2025-12-14 19:41:18)
#    for the models that already existed in the sheet (rows 6-11).
#    Style B2 = green "OK" look, B3 = yellow "429" look.
# ---------------------------------------------------------------------
function Set-Result($cellRef, $styleRef, $value) {
    $ws.Range($styleRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
    $ws.Range($cellRef).Value = $value
}

Set-Result "F6" "B2" "OK (1584ms)"
Set-Result "G6" "B2" "OK (1257ms)"

Set-Result "F7" "B3" "429"
Set-Result "G7" "B3" "429"

Set-Result "F8" "B2" "OK (2903ms)"
Set-Result "G8" "B2" "OK (5822ms)"

Set-Result "F9" "B2" "OK (1472ms)"
Set-Result "G9" "B2" "OK (1501ms)"

Set-Result "F10" "B2" "OK (1663ms)"
Set-Result "G10" "B2" "OK (1286ms)"

Set-Result "F11" "B2" "OK (1242ms)"
Set-Result "G11" "B2" "OK (1282ms)"

# ---------------------------------------------------------------------
# 5. New model rows (12-15) that only appear starting with run F.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A12:A15").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("B12:E15").PasteSpecial(-4122)

$ws.Range("A12").Value = "google/gemma-3n-e4b-it:free"
$ws.Range("A13").Value = "mistralai/devstral-2512:free"
$ws.Range("A14").Value = "mistralai/mistral-7b-instruct:free"
$ws.Range("A15").Value = "mistralai/mistral-small-3.1-24b-instruct:free"

Set-Result "F12" "B2" "OK (1266ms)"
Set-Result "G12" "B2" "OK (1485ms)"

Set-Result "F13" "B2" "OK (3733ms)"
Set-Result "G13" "B3" "429"

Set-Result "F14" "B2" "OK (692ms)"
Set-Result "G14" "B2" "OK (1290ms)"

Set-Result "F15" "B2" "OK (1418ms)"
Set-Result "G15" "B2" "OK (937ms)"

Write-Host "done"
